$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch the price/volume columns to text format so that numeric-looking
# strings (e.g. "0.9978") are stored as text, matching the source data,
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.732.87"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "1.725.54"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("D4").Value = "0.9978"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "240.35"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").Value = "0.9985"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "0.4811"
$ws.Range("E7").Value = "  -1.45%  "
$ws.Range("D8").Value = "0.2586"
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").Value = "0.06179"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "1.722.55"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("D11").Value = "15.81"
$ws.Range("E11").Value = "  +1.99%  "
$ws.Range("D12").Value = "0.06848"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("D13").Value = "0.6022"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").Value = "76.87"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "0.9985"
$ws.Range("D17").Value = "26.553.59"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").Value = "0.9980"
$ws.Range("D19").Value = "0.000007136"
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").Value = "1.944.81"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").Value = "8.511"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").Value = "5.051"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("E27").Value = "  +2.68%  "
$ws.Range("D28").Value = "106.21"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  -2.40%  "
$ws.Range("D30").Value = "4.002"
$ws.Range("E30").Value = "  +2.67%  "
$ws.Range("D31").Value = "0.07905"
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("D32").Value = "3.658"
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").Value = "0.04526"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").Value = "2.596"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").Value = "0.9970"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("D37").Value = "0.9305"
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("D38").Value = "2.452"
$ws.Range("E38").Value = "  +3.05%  "
$ws.Range("D39").Value = "1.988"
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("D40").Value = "0.9976"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("D41").Value = "0.01496"
$ws.Range("E41").Value = "  +1.37%  "
$ws.Range("D42").Value = "5.603"
$ws.Range("E42").Value = "  +3.19%  "
$ws.Range("D43").Value = "99.93"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D44").Value = "0.3821"
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("D45").Value = "6.763"
$ws.Range("E45").Value = "  -1.99%  "
$ws.Range("D46").Value = "0.1152"
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("D47").Value = "0.05356"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").Value = "7.889"
$ws.Range("E48").Value = "  +2.83%  "
$ws.Range("D49").Value = "30.06"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").Value = "1.247"
$ws.Range("E50").Value = "  +2.44%  "
$ws.Range("D51").Value = "51.38"

# Remove the temporary text-number-format override so the cells keep their
# original (default/general) style, same as before the edit.
$ws.Range("D2:E51").ClearFormats()
